# Remove the obsolete "x" marker column (column A) from the routes sheet.
# It was duplicated on almost every data row (rows 5-40) and is no longer
# needed now that the authorController has its own dedicated routes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5:A40").ClearContents()

# Reset the view: select A9 (top of the visible area) instead of the old
# C30 selection further down the sheet.
$ws.Range("A9").Select()
